$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch the lowercase default date/time format first (as Excel itself
# proposes when a date format is first applied), then settle on the
# final, explicit uppercase format for the whole column of dates.
$ws.Range("K2").NumberFormat = "yyyy-mm-dd h:mm:ss"

$rng = $ws.Range("K2:K41")
$rng.Value = 45510
$rng.NumberFormat = "YYYY-MM-DD HH:MM:SS"
